# Regenerate the K column (column G) values to reflect K (strikeouts-like
# count) instead of the old "Strike#" derived values, as described in the
# commit message: "regen save_data to use K instead of Strike#, regen
# std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$updates = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 0
    6  = 1
    9  = 1
    10 = 1
    11 = 2
    12 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
